# Edit script applying the cre.docx diff:
#  1. Split the "Final Certified Application to RIC" paragraph into three
#     paragraphs: the original heading (now bold / navy), a new
#     "Today please" paragraph (which also carries the _GoBack bookmark),
#     and a new "Here's something new that you have to do." paragraph
#     (which keeps the original trailing space run).
#  2. Move a <w:lastRenderedPageBreak/> from the "The intent of the CRE in
#     Infectious Disease..." run to the "CRE in He" run.
#  3. Add a <w:lastRenderedPageBreak/> to the "CIA and Head of
#     Department/Institute Director signatures." run.
#  4. Remove the <w:lastRenderedPageBreak/> before "written evidence must be
#     attached..." and merge that run back into the preceding one.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: split "Final Certified Application to RIC" paragraph into 3
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Final Certified Application to RIC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r.Find.Found) { throw "anchor 1 not found" }
$fullPara = $r.Paragraphs(1).Range

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00000000" w:rsidRDefault="00856ACB"><w:pPr><w:pBdr><w:left w:val="single" w:sz="6" w:space="0" w:color="D5E0E9"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:spacing w:line="360" w:lineRule="atLeast"/><w:ind w:left="720"/><w:divId w:val="258292435"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:b/><w:bCs/><w:color w:val="0C304A"/><w:spacing w:val="3"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:b/><w:bCs/><w:color w:val="0C304A"/><w:spacing w:val="3"/></w:rPr><w:t>Final Certified Application to RIC</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:left w:val="single" w:sz="6" w:space="0" w:color="D5E0E9"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:spacing w:line="360" w:lineRule="atLeast"/><w:ind w:left="720"/><w:divId w:val="258292435"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="666666"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:caps/><w:color w:val="666666"/><w:sz w:val="23"/><w:szCs w:val="23"/></w:rPr><w:t>Today please</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pBdr><w:left w:val="single" w:sz="6" w:space="0" w:color="D5E0E9"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:spacing w:line="360" w:lineRule="atLeast"/><w:ind w:left="720"/><w:divId w:val="258292435"/><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="666666"/><w:spacing w:val="3"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:b/><w:bCs/><w:color w:val="0C304A"/><w:spacing w:val="3"/></w:rPr><w:t>Here&#8217;s something new that you have to do.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="666666"/><w:spacing w:val="3"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@

$fullPara.InsertXML($xml1)

Write-Output "step1 done"

# ---------------------------------------------------------------------
# Step 2: move <w:lastRenderedPageBreak/> from the "The intent of the CRE
# in Infectious Disease..." run onto the "CRE in He" run.
# ---------------------------------------------------------------------
$r2a = $d.Content
$r2a.Find.Execute("CRE in Health Services Research;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r2a.Find.Found) { throw "anchor 2a not found" }
$para2a = $r2a.Paragraphs(1).Range

$xml2a = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00000000" w:rsidRDefault="00856ACB"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:ind w:left="225"/><w:divId w:val="301082192"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:lastRenderedPageBreak/><w:t>CRE in He</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t>alth Services Research;</w:t></w:r></w:p>
'@

$para2a.InsertXML($xml2a)

$r2b = $d.Content
$r2b.Find.Execute("The intent of the CRE in Infectious Disease Emergency Response Research stream", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r2b.Find.Found) { throw "anchor 2b not found" }
$para2b = $r2b.Paragraphs(1).Range

$xml2b = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00000000" w:rsidRDefault="00856ACB"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:ind w:left="225"/><w:divId w:val="301082192"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t>The intent of the CRE in Infectious Disease Emergency Response Research stream is to provide support for a single team of researchers to e</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t xml:space="preserve">stablish a national research capacity that responds to major infectious disease threats and embeds research within the health system. This team will develop evidence that can be implemented in future epidemics, decrease the burden on the health system and </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t xml:space="preserve">improve health outcomes. Research, training and capacity building will be focussed on research activities occurring before, during and after infectious disease emergencies. </w:t></w:r></w:p>
'@

$para2b.InsertXML($xml2b)

Write-Output "step2 done"

# ---------------------------------------------------------------------
# Step 3: add <w:lastRenderedPageBreak/> to the "CIA and Head of
# Department/Institute Director signatures." run.
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("CIA and Head of Department/Institute Director signatures.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r3.Find.Found) { throw "anchor 3 not found" }
$para3 = $r3.Paragraphs(1).Range

$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00000000" w:rsidRDefault="00856ACB"><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:ind w:left="450"/><w:divId w:val="1535583162"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:lastRenderedPageBreak/><w:t>CIA and Head of Department/Institute Director signatures.</w:t></w:r></w:p>
'@

$para3.InsertXML($xml3)

Write-Output "step3 done"

# ---------------------------------------------------------------------
# Step 4: remove <w:lastRenderedPageBreak/> before "written evidence must
# be attached..." and merge its run back with the preceding one.
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Written evidence (e.g", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r4.Find.Found) { throw "anchor 4 not found" }
$para4 = $r4.Paragraphs(1).Range

$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00000000" w:rsidRDefault="00856ACB"><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="0C304A"/><w:ind w:left="450"/><w:divId w:val="1535583162"/><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t>Written evidence (e.g</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t>. email) from all CIs and AIs clearly stating their agreement to be on the application along with the application ID number. All written evidence must be attached to the finished Authority to Submit form with the other signatures. This is an NHMRC requirem</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:spacing w:val="3"/></w:rPr><w:t>ent and the application cannot be submitted to NHMRC until received.</w:t></w:r></w:p>
'@

$para4.InsertXML($xml4)

Write-Output "step4 done"
